# "commit for first push"
#
# - Rename the shared "satN" station-name strings to "sat_N"
#   (column B / "station", rows 2-8 on Sheet1).
# - Move the active cell selection from B2 to B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "sat_1"
$ws.Range("B3").Value = "sat_2"
$ws.Range("B4").Value = "sat_3"
$ws.Range("B5").Value = "sat_4"
$ws.Range("B6").Value = "sat_5"
$ws.Range("B7").Value = "sat_6"
$ws.Range("B8").Value = "sat_7"

$ws.Range("B10").Select() | Out-Null
